# Adding Master Data XLS
# Re-assign the "lang_code" (E) column for the translated reason-list rows
# to their real ISO language codes instead of the placeholder "eng", and
# strip the stray quote marks from two of the Arabic "demographic
# details" descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the extraneous quotation marks around the "all demographic
#     details" Arabic string (row 17) --------------------------------------
$ws.Range("B17").Value = "جميع تفاصيل ديموغرافية هي مطابقة"
$ws.Range("C17").Value = "جميع تفاصيل ديموغرافية هي مطابقة"

# --- Arabic block (rows 11-19): lang_code eng -> ara ---------------------
$ws.Range("E11:E19").Value = "ara"

# --- French block (rows 20-28): lang_code eng -> fra ----------------------
$ws.Range("E20:E28").Value = "fra"

# --- Drop the extraneous quotation marks around the "some demographic
#     details" Arabic string (row 19) --------------------------------------
$ws.Range("B19").Value = "بعض التفاصيل الديمغرافية هي مطابقة"
$ws.Range("C19").Value = "بعض التفاصيل الديمغرافية هي مطابقة"

# --- Cosmetic touch-ups that came along with the data edit ---------------
# Column widths so the long Arabic / French text is fully visible.
$ws.Columns.Item(2).ColumnWidth = 40.1
$ws.Columns.Item(3).ColumnWidth = 44.92

# Scroll the view down and leave the selection on the row that was edited.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select() | Out-Null

# Page setup as left by the author (Print Setup dialog touched once).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
